# Fix the mis-typed "2050" (resp. "2041-2050") column header that was left
# over as a stray numeric value (709.3600745874514) on the E1 header cell
# of each table, and drop the "Total" row that trails each table.

$wb = $excel.ActiveWorkbook

# Sheets 1-3 and 5 use plain year labels (2015 / 2030 / 2040 / 2050).
# The label is purely numeric, so it must be entered with a leading
# apostrophe to force Excel to store it as text instead of a number.
$yearLabelSheets = @(1, 2, 3, 5)
foreach ($idx in $yearLabelSheets) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Range("E1").Value = "'2050"
}

# Sheet 4 uses period labels (2015-2030 / 2031-2040 / 2041-2050); this text
# is not numeric-looking so it can be assigned directly.
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("E1").Value = "2041-2050"

# Remove the trailing "Total" row from each table that has one (sheets
# 1-4 have it at row 13, sheet 6 has it at row 4). Sheet 5 never had a
# Total row.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(13).Delete()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(13).Delete()

$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(13).Delete()

$ws4.Rows.Item(13).Delete()

$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(4).Delete()
